$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 1 - header translations (PT -> ES)
# ---------------------------------------------------------------
$ws.Range("B1").Value = "Dirección"
$ws.Range("C1").Value = "Latitud"
$ws.Range("D1").Value = "Longitud"
$ws.Range("F1").Value = "Barrio"
$ws.Range("G1").Value = "Referencia"
$ws.Range("H1").Value = "Ciudad"
$ws.Range("J1").Value = "Provincia"
$ws.Range("O1").Value = "Valor de la tabla"
$ws.Range("Q1").Value = "Producción"
$ws.Range("R1").Value = "Comentarios"
$ws.Range("U1").Value = "Costo neto"
$ws.Range("W1").Value = "Comentarios internos"

# ---------------------------------------------------------------
# Row 2 - replace data (new location record)
# ---------------------------------------------------------------
$ws.Range("A2").Value = "MDV40"
$ws.Range("B2").Value = "Av. Torquato Tapajós,  3741, Prox VIDEOLAR e AMBEV - Novo Israel - Sent Centro"
$ws.Range("C2").Value = -3.0251296
$ws.Range("D2").Value = -60.0214348
$ws.Range("E2").Value = "https://i.ibb.co/Xk6bLnt/Av-Torquato-Tapaj-s-3741-Prox-VIDEOLAR-e-AMBEV-Novo-Israel.jpg"
$ws.Range("H2").Value = "Manaus"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "10x4"

# ---------------------------------------------------------------
# Row 3 - replace data (new location record)
# ---------------------------------------------------------------
$ws.Range("A3").Value = "MMK10"
$ws.Range("B3").Value = "AV. PERIMETRAL, (APAE)  - SENTIDO BOLDA DO ELDORADO - próximo do Shopping Manauara e em frente ao CSU"
$ws.Range("C3").Value = -3.0857
$ws.Range("D3").Value = -60.01127
$ws.Range("E3").Value = "https://i.ibb.co/bRrJQVp/AV-PERIMETRAL-APAE-SENTIDO-BOLDA-DO-ELDORADO.png"
$ws.Range("F3").Value = ""
$ws.Range("H3").Value = "Manaus"
$ws.Range("J3").Value = ""
$ws.Range("L3").Value = "10x4"
